$d = $word.ActiveDocument

# Locate the "4-payment cancellation ..." paragraph (the paragraph that
# immediately precedes the bold "lookup types:" heading) by searching the
# document's paragraphs for the unique lead-in text.
$targetIndex = -1
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text -like "*payment cancellation happens thru invoice*") {
        $targetIndex = $i
    }
}

$anchor = $d.Paragraphs($targetIndex)
$anchorRange = $anchor.Range

# Insert two new (non-bold) paragraphs right after the anchor paragraph and
# before the "lookup types:" heading.
$anchorRange.InsertParagraphAfter()
$anchorRange.InsertParagraphAfter()

$newPara1 = $d.Paragraphs($targetIndex + 1)
$newPara2 = $d.Paragraphs($targetIndex + 2)

$newPara1.Range.Text = "5-in polymorphism if a method is defined virtual-override, even if method is called in base class, will bubble up to inherited class. which is fine but tricky"
$newPara2.Range.Text = " "

Write-Output "Inserted commentary paragraphs before 'lookup types:' heading."
